$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update flight number and departure time (date serial 45635.375)
$ws.Range("A2").Value = "VN0012"
$ws.Range("C2").Value = 45635.375

# Rows 3-6 held placeholder flights (VN0003..VN0007) that are no longer valid
# data for the import-validation flow; clear them back out, leaving only the
# pre-formatted (date-styled) C column like the empty rows below them.
$ws.Range("A3:I6").ClearContents()

# Move/save the active selection to A2 (was H10)
$ws.Range("A2").Select()
